$d = $word.ActiveDocument

# Replace the date-line text (occurs 4 times): "Pegaz" -> "Pegasus" and hyphens -> en-dashes
$d.Content.Find.Execute(
    "2022: Datumi kampanje za opazovanje ozvezdje Pegaz: 8.-17. oktober, 7.-16. november,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2022: Datumi kampanje za opazovanje ozvezdje Pegasus: 8." + [char]8211 + "17. oktober, 7." + [char]8211 + "16. november,",
    2
)

# Replace the paragraph text (occurs once): "Pegaz" -> "Pegasus"
$d.Content.Find.Execute(
    "Sodelujete v svetovni aktivnosti opazovanja in beleženja najšibkejših, s prostim očesom  še vidnih zvezd, kot metode za merjenje svetlobnega onesnaževanja na določenem mestu. Z opazovanjem izbranega ozvezdje Pegaz na nočnem nebu in s primerjavo videnega z zvezdnimi kartami, se lahko ljudje širom sveta podučijo o tem, kako svetila v njihovem kraju prispevajo k svetlobnemu onesnaževanju.  Vaši prispevki v spletno bazo podatkov bodo pomagali dokumentirati nočno nebo, vidno s prostim očesom.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Sodelujete v svetovni aktivnosti opazovanja in beleženja najšibkejših, s prostim očesom  še vidnih zvezd, kot metode za merjenje svetlobnega onesnaževanja na določenem mestu. Z opazovanjem izbranega ozvezdje Pegasus na nočnem nebu in s primerjavo videnega z zvezdnimi kartami, se lahko ljudje širom sveta podučijo o tem, kako svetila v njihovem kraju prispevajo k svetlobnemu onesnaževanju.  Vaši prispevki v spletno bazo podatkov bodo pomagali dokumentirati nočno nebo, vidno s prostim očesom.",
    2
)
